# Update countries & provincias Spain
# Applies the data refresh captured by the source diff:
#  - Refreshed case counts for several countries (Brasil, Pakistan, Corea del Sur,
#    Costa Rica) and a re-sort around Kirguistan/Malasia/El Salvador/Kenia caused
#    by Kirguistan's case count overtaking the others.
#  - Two tie-break reorderings (Seychelles/Lesoto and Islas Malvinas/Groenlandia)
#    where the underlying numbers are identical but the listed order swapped.
#  - Updated "Datos actualizados" timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: footer timestamp text
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 06:43"

# Row 5: Brasil - Casos activos / Recuperados refreshed
$ws.Range("D5").Value = 1152467
$ws.Range("E5").Value = 495674

# Row 15: Pakistan - full refresh of figures
$ws.Range("B15").Value = 240848
$ws.Range("C15").Value = 3359
$ws.Range("D15").Value = 145311
$ws.Range("E15").Value = 90554
$ws.Range("G15").Value = 61
$ws.Range("H15").Value = 4983

# Row 67: Corea del Sur - Nuevos casos
$ws.Range("C67").Value = 50

# Rows 75-78: Kirguistan jumps ahead of Malasia / El Salvador / Kenia because its
# total cases (8847) now exceed Malasia's (8677), so the block shifts down by one.
$ws.Range("A75").Value = "Kirguistan"
$ws.Range("B75").Value = 8847
$ws.Range("C75").Value = 361
$ws.Range("D75").Value = 3053
$ws.Range("E75").Value = 5678
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 116

$ws.Range("A76").Value = "Malasia"
$ws.Range("B76").Value = 8677
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 8486
$ws.Range("E76").Value = 70
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 121

$ws.Range("A77").Value = "El Salvador"
$ws.Range("B77").Value = 8566
$ws.Range("C77").Value = 0
$ws.Range("D77").Value = 5133
$ws.Range("E77").Value = 3198
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 235

$ws.Range("A78").Value = "Kenia"
$ws.Range("B78").Value = 8528
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 2593
$ws.Range("E78").Value = 5766
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 169

# Row 90: Costa Rica - minor refresh
$ws.Range("E90").Value = 3882
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 25

# Rows 184-185: Seychelles / Lesoto tie-break swap (identical figures, order flips)
$ws.Range("A184").Value = "Seychelles"
$ws.Range("A185").Value = "Lesoto"

# Rows 209-210: Islas Malvinas / Groenlandia tie-break swap (identical figures, order flips)
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
